$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("creditvouchercreatepayment")

# Remove the obsolete "rfid" column (value "Kul-01222") - credit voucher
# payments no longer carry an rfid, shifting refType/currency/amount/... left.
$ws.Columns("B:B").Delete()

# Best-effort cosmetic follow-up resize of the (now shifted) paymentMethodCode
# column, applied below the data rows so the existing header/data cell styles
# stay untouched.
$ws.Range("F3:F10").ColumnWidth = 23

# The credit-voucher tab becomes the active/selected sheet, with the cursor
# left on D19.
$ws.Activate() | Out-Null
$ws.Range("D19").Select() | Out-Null
